$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 485, pushing the existing rows 485-506 down to 486-507.
$ws.Rows.Item(485).Insert()

# Populate the new row 485 with this week's entry. Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T
# carry the same values as the (now shifted) row below it; D,M,N,O,P,S are new.
$ws.Range("A485").Value = 4
$ws.Range("B485").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C485").Value = "Los Lagos"
$ws.Range("D485").Value = 45267
$ws.Range("E485").Value = 10
$ws.Range("F485").Value = "Fruta"
$ws.Range("G485").Value = 100108
$ws.Range("H485").Value = "Tropicales y subtropicales"
$ws.Range("I485").Value = 100108005
$ws.Range("J485").Value = "Piña"
$ws.Range("K485").Value = "Caramelo"
$ws.Range("L485").Value = "Segunda"
$ws.Range("M485").Value = 100
$ws.Range("N485").Value = 28000
$ws.Range("O485").Value = 28000
$ws.Range("P485").Value = 28000
$ws.Range("Q485").Value = "$/caja 14 unidades"
$ws.Range("R485").Value = "Ecuador"
$ws.Range("S485").Value = 2000
$ws.Range("T485").Value = 14
